$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F6/F7: record moved from "En revisión por parte de editor" into the next
# tracking stage -> "En manuscrito de autor" (adds a new shared string).
$ws.Range("F6").Value = "En manuscrito de autor"
$ws.Range("F7").Value = "En manuscrito de autor"

# Row 12 (item 7) dates were entered by mistake - clear them back out while
# keeping the existing cell formatting/style.
$ws.Range("B12:E12").ClearContents()

# Reflect where the coordinator left off reviewing: scroll the view down a
# bit and leave the selection on F7.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F7").Select()
